$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @{
    2 = @(3, 1, 57.52065133333334, 172.561954, 0.1828443315107865, 0.1855832454108249, 3, 1, 0.1255626666666667, 0.376688, 0.02744849445093922, 0.02907345870642374, 7.222446369816891, 65.00201732835201, 0.005018801618859513, 0.00539554682205572)
    3 = @(3, 1, 57.52065133333334, 172.561954, 0.1828443315107865, 0.1855832454108249, 3, 1, 3.393572666666666, 10.180718, 0.7418483772500767, 0.7857661629113346, 195.2005101336636, 1756.804591202972, 0.135642770620652, 0.1458250346470964)
    4 = @(3, 1, 57.52065133333334, 172.561954, 0.1828443315107865, 0.1855832454108249, 3, 1, 0.2723486666666667, 0.817046, 0.05953649332381727, 0.06306108275880487, 15.66567269643156, 140.991054267884, 0.01088591032228977, 0.01170308039749962)
    5 = @(3, 1, 57.52065133333334, 172.561954, 0.1828443315107865, 0.1855832454108249, 2, 1, 0.7670265, 1.534053, 0.167675019875653, 0.1184009752075072, 44.11986386992701, 264.7191832195621, 0.03065842692022161, 0.02197323723881579)
    6 = @(3, 1, 57.52065133333334, 172.561954, 0.1828443315107865, 0.1855832454108249, 1, 0.3333333333333333, 0.01597233333333333, 0.047917, 0.0034916150995138, 0.003698320415929645, 0.9187390166464446, 8.268651149818002, 0.0006384220287635689, 0.0006863463053573353)
    7 = @(3, 1, 243.0020346666667, 729.006104, 0.7724450880589986, 0.7840159175794992, 3, 1, 0.1255626666666667, 0.376688, 0.02744849445093922, 0.02907345870642374, 30.51198347817245, 274.607851303552, 0.02120245471324268, 0.02279405440492649)
    8 = @(3, 1, 243.0020346666667, 729.006104, 0.7724450880589986, 0.7840159175794992, 3, 1, 3.393572666666666, 10.180718, 0.7418483772500767, 0.7857661629113346, 824.6450627891857, 7421.805565102672, 0.5730371350913607, 0.6160531792178522)
    9 = @(3, 1, 243.0020346666667, 729.006104, 0.7724450880589986, 0.7840159175794992, 3, 1, 0.2723486666666667, 0.817046, 0.05953649332381727, 0.06306108275880487, 66.18128013875378, 595.6315212487841, 0.04598867182824001, 0.04944089266270114)
    10 = @(3, 1, 243.0020346666667, 729.006104, 0.7724450880589986, 0.7840159175794992, 2, 1, 0.7670265, 1.534053, 0.167675019875653, 0.1184009752075072, 186.389000143252, 1118.334000859512, 0.1295197454931431, 0.09282824921962125)
    11 = @(3, 1, 243.0020346666667, 729.006104, 0.7724450880589986, 0.7840159175794992, 1, 0.3333333333333333, 0.01597233333333333, 0.047917, 0.0034916150995138, 0.003698320415929645, 3.881309498374223, 34.931785485368, 0.002697080933012066, 0.002899542074398076)
    12 = @(1, 0.3333333333333333, 0.1029616666666667, 0.308885, 0.000327290402255814, 0.0003321930438891683, 3, 1, 0.1255626666666667, 0.376688, 0.02744849445093922, 0.02907345870642374, 0.01292814143111111, 0.11635327288, 0.000008983628790164376, 0.000009658000744072945)
    13 = @(1, 0.3333333333333333, 0.1029616666666667, 0.308885, 0.000327290402255814, 0.0003321930438891683, 3, 1, 3.393572666666666, 10.180718, 0.7418483772500767, 0.7857661629113346, 0.3494078977144444, 3.14467107943, 0.0002427998538030005, 0.0002610260534426284)
    14 = @(1, 0.3333333333333333, 0.1029616666666667, 0.308885, 0.000327290402255814, 0.0003321930438891683, 3, 1, 0.2723486666666667, 0.817046, 0.05953649332381727, 0.06306108275880487, 0.02804147263444445, 0.25237325371, 0.00001948572284885274, 0.00002094845303259414)
    15 = @(1, 0.3333333333333333, 0.1029616666666667, 0.308885, 0.000327290402255814, 0.0003321930438891683, 2, 1, 0.7670265, 1.534053, 0.167675019875653, 0.1184009752075072, 0.07897432681750001, 0.4738459609050001, 0.00005487842470335409, 0.00003933198035362776)
    16 = @(1, 0.3333333333333333, 0.1029616666666667, 0.308885, 0.000327290402255814, 0.0003321930438891683, 1, 0.3333333333333333, 0.01597233333333333, 0.047917, 0.0034916150995138, 0.003698320415929645, 0.001644538060555556, 0.014800842545, 0.000001142772110442346, 0.000001228556316245124)
    17 = @(2, 1, 13.928462, 27.856924, 0.04427523444762439, 0.02995896976851976, 3, 1, 0.1255626666666667, 0.376688, 0.02744849445093922, 0.02907345870642374, 1.748894831285334, 10.493368987712, 0.001215288527049651, 0.0008710108704520564)
    18 = @(2, 1, 13.928462, 27.856924, 0.04427523444762439, 0.02995896976851976, 3, 1, 3.393572666666666, 10.180718, 0.7418483772500767, 0.7857661629113346, 47.26724793190532, 283.603487591432, 0.03284551082733685, 0.02354074471978645)
    19 = @(2, 1, 13.928462, 27.856924, 0.04427523444762439, 0.02995896976851976, 3, 1, 0.2723486666666667, 0.817046, 0.05953649332381727, 0.06306108275880487, 3.793398054417334, 22.760388326504, 0.002635992200101434, 0.001889245071941158)
    20 = @(2, 1, 13.928462, 27.856924, 0.04427523444762439, 0.02995896976851976, 2, 1, 0.7670265, 1.534053, 0.167675019875653, 0.1184009752075072, 10.683499458243, 42.733997832972, 0.007423850816004616, 0.003547171236804965)
    21 = @(2, 1, 13.928462, 27.856924, 0.04427523444762439, 0.02995896976851976, 1, 0.3333333333333333, 0.01597233333333333, 0.047917, 0.0034916150995138, 0.003698320415929645, 0.2224700378846667, 1.334820227308, 0.0001545920771318388, 0.0001107978695351357)
    22 = @(1, 0.3333333333333333, 0.033993, 0.101979, 0.0001080555803345765, 0.0001096741972668582, 3, 1, 0.1255626666666667, 0.376688, 0.02744849445093922, 0.02907345870642374, 0.004268251728000001, 0.038414265552, 0.00000296596299720664, 0.000003188608245398173)
    23 = @(1, 0.3333333333333333, 0.033993, 0.101979, 0.0001080555803345765, 0.0001096741972668582, 3, 1, 3.393572666666666, 10.180718, 0.7418483772500767, 0.7857661629113346, 0.115357715658, 1.038219440922, 0.00008016085692402087, 0.00008617827315675995)
    24 = @(1, 0.3333333333333333, 0.033993, 0.101979, 0.0001080555803345765, 0.0001096741972668582, 3, 1, 0.2723486666666667, 0.817046, 0.05953649332381727, 0.06306108275880487, 0.009257948226000002, 0.083321534034, 0.000006433250337190714, 0.000006916173630350835)
    25 = @(1, 0.3333333333333333, 0.033993, 0.101979, 0.0001080555803345765, 0.0001096741972668582, 2, 1, 0.7670265, 1.534053, 0.167675019875653, 0.1184009752075072, 0.0260735318145, 0.156441190887, 0.00001811822158027533, 0.00001298553191149653)
    26 = @(1, 0.3333333333333333, 0.033993, 0.101979, 0.0001080555803345765, 0.0001096741972668582, 1, 0.3333333333333333, 0.01597233333333333, 0.047917, 0.0034916150995138, 0.003698320415929645, 0.0005429475270000001, 0.004886527743, 0.0000003772884958829337, 0.000000405610322852717)
}

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item([int]$r, 5 + $i).Value = $vals[$i]
    }
}

Write-Output "done"